$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "US4 - Check Referral Status" and "US5 - Reward Details" rows;
# remaining rows below shift up automatically.
$ws.Range("A11:A12").EntireRow.Delete() | Out-Null

# Update the Team Members for the project.
$ws.Range("D3").Value2 = "PRANALI RAUT,PUJA CHANEKAR"

# Clear out the "Prepared By" name and "Date of Review" date, leaving the cells blank.
$ws.Range("B4").ClearContents() | Out-Null
$ws.Range("D4").ClearContents() | Out-Null

# Reset the view: scroll back to the top and select cell B5.
$ws.Range("B5").Select() | Out-Null
